$d = $word.ActiveDocument

# The document currently ends with the paragraph "Üçüncü paragrafı".
# Add a new paragraph right after it containing the italic note about
# the father/son story.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Text = "-Baba ve oğulun hikayesi.Burada bize verilen görev hikayenin bu kısmını tamamlamak."
$newRange.Font.Italic = $true
